$d = $word.ActiveDocument

# Rewrite the notarization instruction per the Aug 18 2023 feedback:
# trims the "along with a notary ... then sign it after you" phrasing down
# to "with a notary public", and drops the trailing "Do not sign your
# forms before they are notarized." sentence.
$old = "Sign and date the letter along with a notary. Write the date at the top and next to your signature. Your notary must see you sign the form, then sign it after you. Do not sign your forms before they are notarized."
$new = "Sign and date the letter with a notary public. Write the date at the top and next to your signature. The notary public must see you sign the form."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# The paragraph used to be followed by a blank "List Paragraph" spacer;
# remove it now that the instruction text is shorter.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*must see you sign the form*") {
        $next = $d.Paragraphs.Item($i + 1)
        $cleanNext = $next.Range.Text -replace "[\r\a]", ""
        if ($cleanNext -eq "" -and $next.Style.NameLocal -eq "List Paragraph") {
            $next.Range.Delete()
        }
        break
    }
}
